{"js": "// Remove the stray italic \"Zacarias\" paragraph that sits directly between the\n// \"ZEC\" (Heading2) book-code paragraph and the following blank/space paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"Zacarias\") {\n    // Found a paragraph whose entire text is \"Zacarias\". Confirm it's the stray\n    // italic one (directly preceded by the \"ZEC\" Heading2 paragraph), not the\n    // later \"Zacarias\" Heading2 section title paragraph.\n    if (i > 0 && paragraphs.items[i - 1].text === \"ZEC\") {\n      target = para;\n      break;\n    }\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the stray italic \"Zacarias\" paragraph that sits directly between the\n# \"ZEC\" (Heading 2) book-code paragraph and the following blank/space paragraph.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$target = $null\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    # Paragraph.Range.Text includes the trailing paragraph-mark control\n    # character(s); strip them off before comparing.\n    $t = $t -replace \"[\\r\\a\\v]+$\", \"\"\n\n    if ($t -eq \"Zacarias\") {\n        # Confirm it's the stray italic paragraph (immediately preceded by the\n        # \"ZEC\" Heading 2 paragraph), not the later \"Zacarias\" Heading 2\n        # section-title paragraph that must stay untouched.\n        if ($i -gt 1) {\n            $prevText = $d.Paragraphs.Item($i - 1).Range.Text -replace \"[\\r\\a\\v]+$\", \"\"\n            if ($prevText -eq \"ZEC\") {\n                $target = $p\n                break\n            }\n        }\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
